$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend the header-row formatting (bold, centered, bordered style used by
# B1:AI1) across the newly added header columns AJ1:BU1 before writing values,
# so new header cells share the same cell style (s="1") as the originals.
$ws.Range("AI1").Copy($ws.Range("AJ1:BU1"))

# Row 1 headers (B1:BU1) -> updated label set
$ws.Range("B1").Value = 'Total Cost'
$ws.Range("C1").Value = 'crudeoil'
$ws.Range("D1").Value = 'natgas'
$ws.Range("E1").Value = 'biomass'
$ws.Range("F1").Value = 'hydrogen'
$ws.Range("G1").Value = 'electricity'
$ws.Range("H1").Value = 'RefineryProduction'
$ws.Range("I1").Value = 'Refinery-gasoline'
$ws.Range("J1").Value = 'Refinery-diesel'
$ws.Range("K1").Value = 'Refinery-kerosene'
$ws.Range("L1").Value = 'Refinery2Production'
$ws.Range("M1").Value = 'Refinery2-gasoline'
$ws.Range("N1").Value = 'Refinery2-diesel'
$ws.Range("O1").Value = 'NGCondProduction'
$ws.Range("P1").Value = 'NGCond-cng'
$ws.Range("Q1").Value = 'BtDProduction'
$ws.Range("R1").Value = 'BtD-diesel'
$ws.Range("S1").Value = 'BtLProduction'
$ws.Range("T1").Value = 'BtL-gasoline'
$ws.Range("U1").Value = 'BM-MethProduction'
$ws.Range("V1").Value = 'BM-Meth-cng'
$ws.Range("W1").Value = 'BM-MeOHProduction'
$ws.Range("X1").Value = 'BM-MeOH-gasoline'
$ws.Range("Y1").Value = 'BM-DMEProduction'
$ws.Range("Z1").Value = 'BM-DME-diesel'
$ws.Range("AA1").Value = 'MtGProduction'
$ws.Range("AB1").Value = 'MtG-gasoline'
$ws.Range("AC1").Value = 'PtF-FT1Production'
$ws.Range("AD1").Value = 'PtF-FT1-gasoline'
$ws.Range("AE1").Value = 'PtF-FT1-diesel'
$ws.Range("AF1").Value = 'PtF-FT2Production'
$ws.Range("AG1").Value = 'PtF-FT2-gasoline'
$ws.Range("AH1").Value = 'PtF-FT2-diesel'
$ws.Range("AI1").Value = 'PtF-FT2-kerosene'
$ws.Range("AJ1").Value = 'PtF-MethProduction'
$ws.Range("AK1").Value = 'PtF-Meth-cng'
$ws.Range("AL1").Value = 'PtF-MeOHProduction'
$ws.Range("AM1").Value = 'PtF-MeOH-gasoline'
$ws.Range("AN1").Value = 'PtF-DMEProduction'
$ws.Range("AO1").Value = 'PtF-DME-diesel'
$ws.Range("AP1").Value = 'PVGasProduction'
$ws.Range("AQ1").Value = 'PVGas-pkm'
$ws.Range("AR1").Value = 'PVDieselProduction'
$ws.Range("AS1").Value = 'PVDiesel-pkm'
$ws.Range("AT1").Value = 'TruckN1Production'
$ws.Range("AU1").Value = 'TruckN1-tkm'
$ws.Range("AV1").Value = 'TruckN2Production'
$ws.Range("AW1").Value = 'TruckN2-tkm'
$ws.Range("AX1").Value = 'TruckN3Production'
$ws.Range("AY1").Value = 'TruckN3-tkm'
$ws.Range("AZ1").Value = 'TruckSZMProduction'
$ws.Range("BA1").Value = 'TruckSZM-tkm'
$ws.Range("BB1").Value = 'PVcngProduction'
$ws.Range("BC1").Value = 'PVcng-pkm'
$ws.Range("BD1").Value = 'TruckN1_ngProduction'
$ws.Range("BE1").Value = 'TruckN1_ng-tkm'
$ws.Range("BF1").Value = 'TruckN2_ngProduction'
$ws.Range("BG1").Value = 'TruckN2_ng-tkm'
$ws.Range("BH1").Value = 'TruckN3_ngProduction'
$ws.Range("BI1").Value = 'TruckN3_ng-tkm'
$ws.Range("BJ1").Value = 'TruckSZM_ngProduction'
$ws.Range("BK1").Value = 'TruckSZM_ng-tkm'
$ws.Range("BL1").Value = 'GasolineHubUsage'
$ws.Range("BM1").Value = 'DieselHubUsage'
$ws.Range("BN1").Value = 'KeroseneHubUsage'
$ws.Range("BO1").Value = 'MethaneHubUsage'
$ws.Range("BP1").Value = 'tkm-N2Usage'
$ws.Range("BQ1").Value = 'tkm-N3Usage'
$ws.Range("BR1").Value = 'tkm-SZMUsage'
$ws.Range("BS1").Value = 'pkmUsage'
$ws.Range("BT1").Value = 'tkm-N1Usage'
$ws.Range("BU1").Value = 'keroseneUsage'

# Row 2 data values (A2:BU2)
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = 0
$ws.Range("C2").Value = 0
$ws.Range("D2").Value = 0
$ws.Range("E2").Value = 0
$ws.Range("F2").Value = 2041.829144604735
$ws.Range("G2").Value = 207.0328368011472
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = 0
$ws.Range("N2").Value = 0
$ws.Range("O2").Value = 0
$ws.Range("P2").Value = 0
$ws.Range("Q2").Value = 0
$ws.Range("R2").Value = 0
$ws.Range("S2").Value = 0
$ws.Range("T2").Value = 0
$ws.Range("U2").Value = 0
$ws.Range("V2").Value = 0
$ws.Range("W2").Value = 0
$ws.Range("X2").Value = 0
$ws.Range("Y2").Value = 0
$ws.Range("Z2").Value = 0
$ws.Range("AA2").Value = 0
$ws.Range("AB2").Value = 0
$ws.Range("AC2").Value = 1491.73788300216
$ws.Range("AD2").Value = 460.9884430110843
$ws.Range("AE2").Value = 1030.749439991076
$ws.Range("AF2").Value = 166.8965517241379
$ws.Range("AG2").Value = 40.88965517241379
$ws.Range("AH2").Value = 5.006896551724139
$ws.Range("AI2").Value = 121
$ws.Range("AJ2").Value = 0
$ws.Range("AK2").Value = 0
$ws.Range("AL2").Value = 0
$ws.Range("AM2").Value = 0
$ws.Range("AN2").Value = 50
$ws.Range("AO2").Value = 50
$ws.Range("AP2").Value = 396.2195511974985
$ws.Range("AQ2").Value = 396.2195511974985
$ws.Range("AR2").Value = 442.5804488025015
$ws.Range("AS2").Value = 442.5804488025015
$ws.Range("AT2").Value = 8
$ws.Range("AU2").Value = 8
$ws.Range("AV2").Value = 26.2
$ws.Range("AW2").Value = 26.2
$ws.Range("AX2").Value = 123.8
$ws.Range("AY2").Value = 123.8
$ws.Range("AZ2").Value = 388.5
$ws.Range("BA2").Value = 388.5
$ws.Range("BB2").Value = 0
$ws.Range("BC2").Value = 0
$ws.Range("BD2").Value = 0
$ws.Range("BE2").Value = 0
$ws.Range("BF2").Value = 0
$ws.Range("BG2").Value = 0
$ws.Range("BH2").Value = 0
$ws.Range("BI2").Value = 0
$ws.Range("BJ2").Value = 0
$ws.Range("BK2").Value = 0
$ws.Range("BL2").Value = 501.8780981834981
$ws.Range("BM2").Value = 1085.7563365428
$ws.Range("BN2").Value = 121
$ws.Range("BO2").Value = 0
$ws.Range("BP2").Value = 26.2
$ws.Range("BQ2").Value = 123.8
$ws.Range("BR2").Value = 388.5
$ws.Range("BS2").Value = 838.8
$ws.Range("BT2").Value = 8
$ws.Range("BU2").Value = 121
